$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds attendance dates as plain text (e.g. "28/07/2022").
# The tutorial output now reports the dates using hyphens instead of
# slashes ("28-07-2022"). Because several of these hyphenated strings
# are ambiguous day/month values, Excel would normally auto-convert
# them into real dates when a value is typed in. Temporarily marking
# the range as Text keeps them as plain strings, matching the source
# data, and the format is reset back to Normal afterwards so no visible
# formatting change is left behind.
$dateRange = $ws.Range("A3:A21")
$dateRange.NumberFormat = "@"

$ws.Range("A3").Value  = "28-07-2022"
$ws.Range("A4").Value  = "01-08-2022"
$ws.Range("A5").Value  = "04-08-2022"
$ws.Range("A6").Value  = "08-08-2022"
$ws.Range("A7").Value  = "11-08-2022"
$ws.Range("A8").Value  = "15-08-2022"
$ws.Range("A9").Value  = "18-08-2022"
$ws.Range("A10").Value = "22-08-2022"
$ws.Range("A11").Value = "25-08-2022"
$ws.Range("A12").Value = "29-08-2022"
$ws.Range("A13").Value = "01-09-2022"
$ws.Range("A14").Value = "05-09-2022"
$ws.Range("A15").Value = "08-09-2022"
$ws.Range("A16").Value = "12-09-2022"
$ws.Range("A17").Value = "15-09-2022"
$ws.Range("A18").Value = "19-09-2022"
$ws.Range("A19").Value = "22-09-2022"
$ws.Range("A20").Value = "26-09-2022"
$ws.Range("A21").Value = "29-09-2022"

$dateRange.Style = "Normal"

# Recomputed attendance counters for the affected sessions.
$ws.Range("D3").Value  = 1
$ws.Range("G3").Value  = 1

$ws.Range("D4").Value  = 1
$ws.Range("E4").Value  = 1
$ws.Range("H4").Value  = 0

$ws.Range("D5").Value  = 1
$ws.Range("E5").Value  = 1
$ws.Range("H5").Value  = 0

$ws.Range("D6").Value  = 1
$ws.Range("G6").Value  = 1

$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("H20").Value = 0
